$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Ligand/Receptor expressing-cell counts (E,K) and recomputed dependent
# expression/specificity/edge-weight statistics for the Fn1-Cd44 LR pair,
# per the NATMI re-run referenced in the commit message.
$rowData = @{
    2 = @{ "E"=2; "G"=40.7349555; "H"=81.469911; "I"=0.05567871843833241; "J"=0.03826666865920979; "K"=2; "M"=19.5719925; "N"=39.143985; "O"=0.03094210933382397; "P"=0.02187190777676379; "Q"=797.2642435338337; "R"=3189.056974135335; "S"=0.001722816993486082; "T"=0.0008369650478382137 }
    3 = @{ "E"=2; "G"=40.7349555; "H"=81.469911; "I"=0.05567871843833241; "J"=0.03826666865920979; "K"=3; "M"=115.495743; "N"=346.487229; "O"=0.1825916246134488; "P"=0.1936015640337701; "Q"=4704.713951544437; "R"=28228.28370926662; "S"=0.0101664676560499; "T"=0.007408486902785067 }
    4 = @{ "E"=2; "G"=40.7349555; "H"=81.469911; "I"=0.05567871843833241; "J"=0.03826666865920979; "K"=3; "M"=239.8982746666667; "N"=719.694824; "O"=0.3792643310961689; "P"=0.4021332732970914; "Q"=9772.245543073444; "R"=58633.47325844067; "S"=0.02111695190480607; "T"=0.01538830072610325 }
    5 = @{ "E"=2; "G"=40.7349555; "H"=81.469911; "I"=0.05567871843833241; "J"=0.03826666865920979; "K"=3; "M"=151.102183; "N"=453.306549; "O"=0.2388832034840335; "P"=0.2532874216646837; "Q"=6155.140700457856; "R"=36930.84420274713; "S"=0.01330071062643437; "T"=0.009692465840388007 }
    6 = @{ "E"=2; "G"=40.7349555; "H"=81.469911; "I"=0.05567871843833241; "J"=0.03826666865920979; "K"=3; "M"=18.12446233333333; "N"=54.373387; "O"=0.02865365369084289; "P"=0.03038141635232813; "Q"=738.2991666097595; "R"=4429.794999658557; "S"=0.001595398716081926; "T"=0.001162595592952039 }
    7 = @{ "E"=2; "G"=40.7349555; "H"=81.469911; "I"=0.05567871843833241; "J"=0.03826666865920979; "K"=2; "M"=88.3431645; "N"=176.686329; "O"=0.139665077781682; "P"=0.09872441687536272; "Q"=3598.65487463668; "R"=14394.61949854672; "S"=0.007776372541474066; "T"=0.003777854549143205 }
    8 = @{ "E"=3; "G"=350.3919066666667; "H"=1051.17572; "I"=0.4789344206933965; "J"=0.4937404802104949; "K"=2; "M"=19.5719925; "N"=39.143985; "O"=0.03094210933382397; "P"=0.02187190777676379; "Q"=6857.8677693407; "R"=41147.2066160442; "S"=0.01481924120882672; "T"=0.01079904624881901 }
    9 = @{ "E"=3; "G"=350.3919066666667; "H"=1051.17572; "I"=0.4789344206933965; "J"=0.4937404802104949; "K"=3; "M"=115.495743; "N"=346.487229; "O"=0.1825916246134488; "P"=0.1936015640337701; "Q"=40468.77360165332; "R"=364218.9624148799; "S"=0.08744941395770821; "T"=0.09558892919553653 }
    10 = @{ "E"=3; "G"=350.3919066666667; "H"=1051.17572; "I"=0.4789344206933965; "J"=0.4937404802104949; "K"=3; "M"=239.8982746666667; "N"=719.694824; "O"=0.3792643310961689; "P"=0.4021332732970914; "Q"=84058.41386649704; "R"=756525.7247984733; "S"=0.1816427427032122; "T"=0.1985494754663241 }
    11 = @{ "E"=3; "G"=350.3919066666667; "H"=1051.17572; "I"=0.4789344206933965; "J"=0.4937404802104949; "K"=3; "M"=151.102183; "N"=453.306549; "O"=0.2388832034840335; "P"=0.2532874216646837; "Q"=52944.98200286559; "R"=476504.8380257902; "S"=0.1144093886740083; "T"=0.125058253203999 }
    12 = @{ "E"=3; "G"=350.3919066666667; "H"=1051.17572; "I"=0.4789344206933965; "J"=0.4937404802104949; "K"=3; "M"=18.12446233333333; "N"=54.373387; "O"=0.02865365369084289; "P"=0.03038141635232813; "Q"=6350.664914284849; "R"=57155.98422856364; "S"=0.01372322103117304; "T"=0.01500053509927347 }
    13 = @{ "E"=3; "G"=350.3919066666667; "H"=1051.17572; "I"=0.4789344206933965; "J"=0.4937404802104949; "K"=2; "M"=88.3431645; "N"=176.686329; "O"=0.139665077781682; "P"=0.09872441687536272; "Q"=30954.72985012198; "R"=185728.3791007319; "S"=0.06689041311846801; "T"=0.04874424099654268 }
    14 = @{ "E"=3; "G"=243.8287033333334; "H"=731.4861100000001; "I"=0.3332781281688242; "J"=0.3435812836494235; "K"=2; "M"=19.5719925; "N"=39.143985; "O"=0.03094210933382397; "P"=0.02187190777676379; "Q"=4772.213552924725; "R"=28633.28131754835; "S"=0.01031232828037196; "T"=0.00751477814980231 }
    15 = @{ "E"=3; "G"=243.8287033333334; "H"=731.4861100000001; "I"=0.3332781281688242; "J"=0.3435812836494235; "K"=3; "M"=115.495743; "N"=346.487229; "O"=0.1825916246134488; "P"=0.1936015640337701; "Q"=28161.17725620991; "R"=253450.5953058892; "S"=0.06085379487047482; "T"=0.06651787388725879 }
    16 = @{ "E"=3; "G"=243.8287033333334; "H"=731.4861100000001; "I"=0.3332781281688242; "J"=0.3435812836494235; "K"=3; "M"=239.8982746666667; "N"=719.694824; "O"=0.3792643310961689; "P"=0.4021332732970914; "Q"=58494.08524387719; "R"=526446.7671948947; "S"=0.1264005063489324; "T"=0.1381654662375591 }
    17 = @{ "E"=3; "G"=243.8287033333334; "H"=731.4861100000001; "I"=0.3332781281688242; "J"=0.3435812836494235; "K"=3; "M"=151.102183; "N"=453.306549; "O"=0.2388832034840335; "P"=0.2532874216646837; "Q"=36843.04935172605; "R"=331587.4441655344; "S"=0.07961454690813104; "T"=0.08702481746780483 }
    18 = @{ "E"=3; "G"=243.8287033333334; "H"=731.4861100000001; "I"=0.3332781281688242; "J"=0.3435812836494235; "K"=3; "M"=18.12446233333333; "N"=54.373387; "O"=0.02865365369084289; "P"=0.03038141635232813; "Q"=4419.264149350508; "R"=39773.37734415457; "S"=0.00954963606728184; "T"=0.01043848602942048 }
    19 = @{ "E"=3; "G"=243.8287033333334; "H"=731.4861100000001; "I"=0.3332781281688242; "J"=0.3435812836494235; "K"=2; "M"=88.3431645; "N"=176.686329; "O"=0.139665077781682; "P"=0.09872441687536272; "Q"=21540.59924839837; "R"=129243.5954903902; "S"=0.04654731569363221; "T"=0.03391986187757793 }
    20 = @{ "E"=3; "G"=48.737294; "H"=146.211882; "I"=0.06661674320651284; "J"=0.06867617226847689; "K"=2; "M"=19.5719925; "N"=39.143985; "O"=0.03094210933382397; "P"=0.02187190777676379; "Q"=953.885952638295; "R"=5723.31571582977; "S"=0.002061262551759196; "T"=0.001502078906317269 }
    21 = @{ "E"=3; "G"=48.737294; "H"=146.211882; "I"=0.06661674320651284; "J"=0.06867617226847689; "K"=3; "M"=115.495743; "N"=346.487229; "O"=0.1825916246134488; "P"=0.1936015640337701; "Q"=5628.949982339442; "R"=50660.54984105498; "S"=0.01216365936853411; "T"=0.01329581436302976 }
    22 = @{ "E"=3; "G"=48.737294; "H"=146.211882; "I"=0.06661674320651284; "J"=0.06867617226847689; "K"=3; "M"=239.8982746666667; "N"=719.694824; "O"=0.3792643310961689; "P"=0.4021332732970914; "Q"=11691.99274252208; "R"=105227.9346826988; "S"=0.02526535455202334; "T"=0.02761697395183755 }
    23 = @{ "E"=3; "G"=48.737294; "H"=146.211882; "I"=0.06661674320651284; "J"=0.06867617226847689; "K"=3; "M"=151.102183; "N"=453.306549; "O"=0.2388832034840335; "P"=0.2532874216646837; "Q"=7364.311516912801; "R"=66278.80365221521; "S"=0.01591362102284501; "T"=0.01739481060368216 }
    24 = @{ "E"=3; "G"=48.737294; "H"=146.211882; "I"=0.06661674320651284; "J"=0.06867617226847689; "K"=3; "M"=18.12446233333333; "N"=54.373387; "O"=0.02865365369084289; "P"=0.03038141635232813; "Q"=883.3372493315926; "R"=7950.035243984335; "S"=0.001908813089851229; "T"=0.002086479383172807 }
    25 = @{ "E"=3; "G"=48.737294; "H"=146.211882; "I"=0.06661674320651284; "J"=0.06867617226847689; "K"=2; "M"=88.3431645; "N"=176.686329; "O"=0.139665077781682; "P"=0.09872441687536272; "Q"=4305.606781126863; "R"=25833.64068676118; "S"=0.009304032621499949; "T"=0.006780015060437337 }
    26 = @{ "E"=3; "G"=22.832077; "H"=68.49623099999999; "I"=0.03120810544755168; "J"=0.03217289111905; "K"=2; "M"=19.5719925; "N"=39.143985; "O"=0.03094210933382397; "P"=0.02187190777676379; "Q"=446.8692398034225; "R"=2681.215438820535; "S"=0.0009656446108596517; "T"=0.0007036825074677244 }
    27 = @{ "E"=3; "G"=22.832077; "H"=68.49623099999999; "I"=0.03120810544755168; "J"=0.03217289111905; "K"=3; "M"=115.495743; "N"=346.487229; "O"=0.1825916246134488; "P"=0.1936015640337701; "Q"=2637.007697348211; "R"=23733.0692761339; "S"=0.005698338674776282; "T"=0.006228722040136272 }
    28 = @{ "E"=3; "G"=22.832077; "H"=68.49623099999999; "I"=0.03120810544755168; "J"=0.03217289111905; "K"=3; "M"=239.8982746666667; "N"=719.694824; "O"=0.3792643310961689; "P"=0.4021332732970914; "Q"=5477.375879356483; "R"=49296.38291420834; "S"=0.01183612123734439; "T"=0.0129377900171345 }
    29 = @{ "E"=3; "G"=22.832077; "H"=68.49623099999999; "I"=0.03120810544755168; "J"=0.03217289111905; "K"=3; "M"=151.102183; "N"=453.306549; "O"=0.2388832034840335; "P"=0.2532874216646837; "Q"=3449.976677124091; "R"=31049.79009411682; "S"=0.007455092203978663; "T"=0.008148988639042775 }
    30 = @{ "E"=3; "G"=22.832077; "H"=68.49623099999999; "I"=0.03120810544755168; "J"=0.03217289111905; "K"=3; "M"=18.12446233333333; "N"=54.373387; "O"=0.02865365369084289; "P"=0.03038141635232813; "Q"=413.8191195782663; "R"=3724.372076204397; "S"=0.0008942262458414533; "T"=0.0009774580003459781 }
    31 = @{ "E"=3; "G"=22.832077; "H"=68.49623099999999; "I"=0.03120810544755168; "J"=0.03217289111905; "K"=2; "M"=88.3431645; "N"=176.686329; "O"=0.139665077781682; "P"=0.09872441687536272; "Q"=2017.057934287666; "R"=12102.347605726; "S"=0.004358682474751238; "T"=0.003176249914922747 }
    32 = @{ "E"=2; "G"=25.082339; "H"=50.164678; "I"=0.03428388404538221; "J"=0.02356250409334498; "K"=2; "M"=19.5719925; "N"=39.143985; "O"=0.03094210933382397; "P"=0.02187190777676379; "Q"=490.9113507904575; "R"=1963.64540316183; "S"=0.00106081568852036; "T"=0.0005153569165192606 }
    33 = @{ "E"=2; "G"=25.082339; "H"=50.164678; "I"=0.03428388404538221; "J"=0.02356250409334498; "K"=3; "M"=115.495743; "N"=346.487229; "O"=0.1825916246134488; "P"=0.1936015640337701; "Q"=2896.903378982877; "R"=17381.42027389726; "S"=0.006259950085905434; "T"=0.004561737645023697 }
    34 = @{ "E"=2; "G"=25.082339; "H"=50.164678; "I"=0.03428388404538221; "J"=0.02356250409334498; "K"=3; "M"=239.8982746666667; "N"=719.694824; "O"=0.3792643310961689; "P"=0.4021332732970914; "Q"=6017.209850704446; "R"=36103.25910422667; "S"=0.0130026543498505; "T"=0.00947526689813293 }
    35 = @{ "E"=2; "G"=25.082339; "H"=50.164678; "I"=0.03428388404538221; "J"=0.02356250409334498; "K"=3; "M"=151.102183; "N"=453.306549; "O"=0.2388832034840335; "P"=0.2532874216646837; "Q"=3789.996177646037; "R"=22739.97706587622; "S"=0.008189844048636048; "T"=0.005968085909766905 }
    36 = @{ "E"=2; "G"=25.082339; "H"=50.164678; "I"=0.03428388404538221; "J"=0.02356250409334498; "K"=3; "M"=18.12446233333333; "N"=54.373387; "O"=0.02865365369084289; "P"=0.03038141635232813; "Q"=454.6039084373977; "R"=2727.623450624386; "S"=0.0009823585406133956; "T"=0.0007158622471633495 }
    37 = @{ "E"=2; "G"=25.082339; "H"=50.164678; "I"=0.03428388404538221; "J"=0.02356250409334498; "K"=2; "M"=88.3431645; "N"=176.686329; "O"=0.139665077781682; "P"=0.09872441687536272; "Q"=2215.853200321766; "R"=8863.412801287062; "S"=0.004788261331856472; "T"=0.00232619447673883 }
}

foreach ($r in $rowData.Keys) {
    $cols = $rowData[$r]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$r").Value = $cols[$col]
    }
}